$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (CON/STR category widths) - update B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 ("CON") - remove old D2 value, add new B2/C2 values
$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = 39.663817176861357
$ws.Range("C2").Value = 19.605916062509721

# Row 3 ("STR") - remove old B3 value, update C3 value
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 38.849098010785461

# Restore the selection to match the new, smaller data range
$ws.Range("B1:E3").Select()
